$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 317.92307
$ws.Range("I55").Value = 280.33334
$ws.Range("J55").Value = 402.5
$ws.Range("K55").Value = 280.33334
$ws.Range("L55").Value = 402.5
$ws.Range("M55").Value = -66.33334000000002
$ws.Range("N55").Value = -830.5

$ws.Range("H58").Value = 524.625
$ws.Range("I58").Value = 278.14285
$ws.Range("K58").Value = 834.4285500000001
$ws.Range("M58").Value = -684.4285500000001

$ws.Range("H61").Value = 846.625
$ws.Range("I61").Value = 846.625
$ws.Range("K61").Value = 2539.875
$ws.Range("M61").Value = -2367.875

$ws.Range("H74").Value = 6927.52
$ws.Range("I74").Value = 3258
$ws.Range("K74").Value = 3258
$ws.Range("M74").Value = -2322

$ws.Range("H77").Value = 6927.52
$ws.Range("I77").Value = 3258
$ws.Range("K77").Value = 16290
$ws.Range("M77").Value = -11610

$ws.Range("H80").Value = 8334900
$ws.Range("I80").Value = 15625862
$ws.Range("K80").Value = 46877586
$ws.Range("M80").Value = -46876588

$ws.Range("H83").Value = 8334900
$ws.Range("I83").Value = 15625862
$ws.Range("K83").Value = 140632758
$ws.Range("M83").Value = -140627766

$ws.Range("H112").Value = 7402.467
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 7402.467
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 22207.401
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -24423.401

$ws.Range("H113").Value = 3649.6428
$ws.Range("I113").Value = 3081.8
$ws.Range("J113").Value = 3965.111
$ws.Range("K113").Value = 3081.8
$ws.Range("L113").Value = 3965.111
$ws.Range("M113").Value = 172.1999999999998
$ws.Range("N113").Value = -10473.111

$ws.Range("H138").Value = 1715.258
$ws.Range("J138").Value = 3550.7144
$ws.Range("L138").Value = 10652.1432
$ws.Range("N138").Value = -20932.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7921.409
$ws.Range("I2").Value = 1787.2727
$ws.Range("J2").Value = 14055.546
$ws.Range("K2").Value = 1787.2727
$ws.Range("L2").Value = 14055.546
$ws.Range("M2").Value = -1674.2727
$ws.Range("N2").Value = -14281.546

$ws.Range("H74").Value = 6680.1113
$ws.Range("I74").Value = 6289.227
$ws.Range("J74").Value = 8400
$ws.Range("K74").Value = 6289.227
$ws.Range("L74").Value = 8400
$ws.Range("M74").Value = -5415.227
$ws.Range("N74").Value = -10148

$ws.Range("H77").Value = 6680.1113
$ws.Range("I77").Value = 6289.227
$ws.Range("J77").Value = 8400
$ws.Range("K77").Value = 31446.135
$ws.Range("L77").Value = 42000
$ws.Range("M77").Value = -27078.135
$ws.Range("N77").Value = -50736

$ws.Range("H116").Value = 7921.409
$ws.Range("I116").Value = 1787.2727
$ws.Range("J116").Value = 14055.546
$ws.Range("K116").Value = 1787.2727
$ws.Range("L116").Value = 14055.546
$ws.Range("M116").Value = 506.7273
$ws.Range("N116").Value = -18643.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7921.409
$ws.Range("I3").Value = 1787.2727
$ws.Range("J3").Value = 14055.546
$ws.Range("K3").Value = 1787.2727
$ws.Range("L3").Value = 14055.546
$ws.Range("M3").Value = -1673.2727
$ws.Range("N3").Value = -14283.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2999.3333
$ws.Range("I62").Value = 2800
$ws.Range("J62").Value = 3099
$ws.Range("K62").Value = 2800
$ws.Range("L62").Value = 3099
$ws.Range("M62").Value = -2176
$ws.Range("N62").Value = -4347

$ws.Range("H65").Value = 2999.3333
$ws.Range("I65").Value = 2800
$ws.Range("J65").Value = 3099
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 15495
$ws.Range("M65").Value = -10880
$ws.Range("N65").Value = -21735

$ws.Range("H86").Value = 83952.32000000001
$ws.Range("I86").Value = 3592.7856
$ws.Range("K86").Value = 3592.7856
$ws.Range("M86").Value = -2469.7856

$ws.Range("H89").Value = 83952.32000000001
$ws.Range("I89").Value = 3592.7856
$ws.Range("K89").Value = 17963.928
$ws.Range("M89").Value = -12347.928

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 51
$ws.Range("I11").Value = 50
$ws.Range("J11").Value = 52
$ws.Range("K11").Value = 150
$ws.Range("L11").Value = 156
$ws.Range("M11").Value = -10
$ws.Range("N11").Value = -436

$ws.Range("H12").Value = 69.5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 69.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 208.5
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -554.5

$ws.Range("H140").Value = 2699.3704
$ws.Range("I140").Value = 1887.2354
$ws.Range("K140").Value = 5661.706200000001
$ws.Range("M140").Value = -481.7062000000005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13496.538
$ws.Range("I80").Value = 16044.3
$ws.Range("J80").Value = 5004
$ws.Range("K80").Value = 16044.3
$ws.Range("L80").Value = 5004
$ws.Range("M80").Value = -15046.3
$ws.Range("N80").Value = -7000

$ws.Range("H83").Value = 13496.538
$ws.Range("I83").Value = 16044.3
$ws.Range("J83").Value = 5004
$ws.Range("K83").Value = 80221.5
$ws.Range("L83").Value = 25020
$ws.Range("M83").Value = -75229.5
$ws.Range("N83").Value = -35004

$ws.Range("H113").Value = 5774.8335
$ws.Range("I113").Value = 2999
$ws.Range("J113").Value = 6330
$ws.Range("K113").Value = 2999
$ws.Range("L113").Value = 6330
$ws.Range("M113").Value = -829
$ws.Range("N113").Value = -10670

$ws.Range("H126").Value = 1044807.1
$ws.Range("I126").Value = 2383330.8
$ws.Range("K126").Value = 7149992.399999999
$ws.Range("M126").Value = -7147522.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1399.6666
$ws.Range("I22").Value = 599.5
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 599.5
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -304.5
$ws.Range("N22").Value = -3590

$ws.Range("H27").Value = 1399.6666
$ws.Range("I27").Value = 599.5
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 599.5
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -492.5
$ws.Range("N27").Value = -3214

$ws.Range("H40").Value = 5852.1816
$ws.Range("I40").Value = 6074.8887
$ws.Range("K40").Value = 6074.8887
$ws.Range("M40").Value = -5938.8887

$ws.Range("H46").Value = 1063
$ws.Range("J46").Value = 1500
$ws.Range("L46").Value = 1500
$ws.Range("N46").Value = -1876

$ws.Range("H55").Value = 210
$ws.Range("I55").Value = 158.33333
$ws.Range("K55").Value = 158.33333
$ws.Range("M55").Value = 14.66667000000001

$ws.Range("H82").Value = 1859.2
$ws.Range("I82").Value = 1924
$ws.Range("J82").Value = 1600
$ws.Range("K82").Value = 1924
$ws.Range("L82").Value = 1600
$ws.Range("M82").Value = -1563
$ws.Range("N82").Value = -2322

$ws.Range("H85").Value = 1859.2
$ws.Range("I85").Value = 1924
$ws.Range("J85").Value = 1600
$ws.Range("K85").Value = 1924
$ws.Range("L85").Value = 1600
$ws.Range("M85").Value = -676
$ws.Range("N85").Value = -4096

$ws.Range("H135").Value = 65000
$ws.Range("J135").Value = 65000
$ws.Range("L135").Value = 65000
$ws.Range("N135").Value = -75140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1662.25
$ws.Range("I81").Value = 1662.25
$ws.Range("K81").Value = 3324.5
$ws.Range("M81").Value = -2263.5

$ws.Range("H84").Value = 1662.25
$ws.Range("I84").Value = 1662.25
$ws.Range("K84").Value = 16622.5
$ws.Range("M84").Value = -11318.5

$ws.Range("H107").Value = 1519.1794
$ws.Range("I107").Value = 978.5333000000001
$ws.Range("J107").Value = 3321.3333
$ws.Range("K107").Value = 2935.5999
$ws.Range("L107").Value = 9963.999899999999
$ws.Range("M107").Value = -1015.5999
$ws.Range("N107").Value = -13803.9999

$ws.Range("H136").Value = 7733.1934
$ws.Range("I136").Value = 7507.8965
$ws.Range("K136").Value = 22523.6895
$ws.Range("M136").Value = -19973.6895

Write-Output "Updated market-price / leve-profit figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
